# Generate Report for Handback
# - Status for the "ea768b0e-..." file moves from "Ready for handoff" to
#   "Handback transform failed" (Overview sheet Status columns + the
#   per-locale "Status" column on the zh-cn / de-de sheets all share the
#   same underlying text).
# - The zh-cn and de-de locale sheets get a new "Error Detail" message
#   explaining the handback/handoff filename mismatch.
# - The "Error Detail" column is widened to fit the new message.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

$newStatus = "Handback transform failed"

# Overview sheet: row 3 is the ea768b0e-d203-40e5-8f4a-728dc4650792.md file;
# columns E (zh-cn) and F (de-de) both show the aggregated status.
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# zh-cn / de-de sheets: column C is "Status" for the same row.
$wsZh.Range("C3").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

# New Error Detail text (column P) for the same row on each locale sheet.
$wsZh.Range("P3").Value = "Handback file name: m5oxjt4u.ivl is different with handoff file name: ea768b0e-d203-40e5-8f4a-728dc4650792.b7305f4becce53e4c20aa55d8d390dc671f9d62a.zh-cn."
$wsDe.Range("P3").Value = "Handback file name: m5oxjt4u.ivl is different with handoff file name: ea768b0e-d203-40e5-8f4a-728dc4650792.b7305f4becce53e4c20aa55d8d390dc671f9d62a.de-de."

# Widen the "Error Detail" column (P) on both locale sheets so the new
# message is readable. 39.17 is the ColumnWidth value that round-trips to
# a stored column width of 40 (matching the other wide columns already on
# this sheet, e.g. column A/G/I/J).
$wsZh.Columns.Item(16).ColumnWidth = 39.17
$wsDe.Columns.Item(16).ColumnWidth = 39.17
